# Minor fix to static decisions_made variable (and related run stats)
# on the "Knowledge Based Agent" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knowledge Based Agent")

# Row 3 (Run #1)
$ws.Range("B3").Value = 21411
$ws.Range("C3").Value = 0
$ws.Range("G3").Value = 16
$ws.Range("H3").Value = -54

# Row 4 (Run #2)
$ws.Range("B4").Value = 146
$ws.Range("C4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# Row 5 (Run #3)
$ws.Range("B5").Value = 768
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = -7

# Row 6 (Run #4)
$ws.Range("B6").Value = 1617
$ws.Range("H6").Value = 989

# Row 7 (Run #5)
$ws.Range("B7").Value = 18292
$ws.Range("D7").Value = 0
$ws.Range("G7").Value = 19
$ws.Range("H7").Value = 954

# Row 8 (Run #6)
$ws.Range("B8").Value = 10131
$ws.Range("D8").Value = 0
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = 976

# Row 9 (Run #7)
$ws.Range("B9").Value = 768
$ws.Range("C9").Value = 0
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = -7

# Row 10 (Run #8)
$ws.Range("B10").Value = 17114
$ws.Range("G10").Value = 19
$ws.Range("H10").Value = 960

# Row 11 (Run #9)
$ws.Range("B11").Value = 10871
$ws.Range("D11").Value = 1
$ws.Range("G11").Value = 24
$ws.Range("H11").Value = 954

# Row 12 (Run #10)
$ws.Range("B12").Value = 6529
$ws.Range("C12").Value = 1
$ws.Range("G12").Value = 17
$ws.Range("H12").Value = 968

# Row 13 (Run #11)
$ws.Range("B13").Value = 2301
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 2
$ws.Range("G13").Value = 17
$ws.Range("H13").Value = 991

# Row 14 (Run #12)
$ws.Range("B14").Value = 8342
$ws.Range("D14").Value = 0
$ws.Range("G14").Value = 13
$ws.Range("H14").Value = 970

# Row 15 (Run #13)
$ws.Range("B15").Value = 5832
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 977

# Row 16 (Run #14)
$ws.Range("B16").Value = 146
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0

# Row 17 (Run #15)
$ws.Range("B17").Value = 24083
$ws.Range("C17").Value = 1
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 944

# Row 18 (Run #16)
$ws.Range("B18").Value = 9426
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 958

# Row 19 (Run #17)
$ws.Range("B19").Value = 15410
$ws.Range("C19").Value = 1
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 953

# Row 20 (Run #18)
$ws.Range("B20").Value = 4716
$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 980

# Row 21 (Run #19)
$ws.Range("B21").Value = 8061
$ws.Range("D21").Value = 1
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 957

# Row 22 (Run #20)
$ws.Range("B22").Value = 5430
$ws.Range("D22").Value = 0
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 981

# Update the view state: clear the scrolled-down "topLeftCell" and move the
# active selection to F13 on the (already) active "Knowledge Based Agent" tab.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F13").Select() | Out-Null
